$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2034.8182
$ws.Range("I116").Value = 1617.7778
$ws.Range("K116").Value = 1617.7778
$ws.Range("M116").Value = 1824.2222
$ws.Range("H132").Value = 2803414
$ws.Range("I132").Value = 2917634.8
$ws.Range("K132").Value = 8752904.399999999
$ws.Range("M132").Value = -8750374.399999999
$ws.Range("H137").Value = 4652284
$ws.Range("I137").Value = 1134.3889
$ws.Range("J137").Value = 8001112
$ws.Range("K137").Value = 3403.1667
$ws.Range("L137").Value = 24003336
$ws.Range("M137").Value = -853.1666999999998
$ws.Range("N137").Value = -24008436
$ws.Range("H138").Value = 4383.0713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29195.785
$ws.Range("I32").Value = 36039.727
$ws.Range("J32").Value = 4101.3335
$ws.Range("K32").Value = 36039.727
$ws.Range("L32").Value = 4101.3335
$ws.Range("M32").Value = -35752.727
$ws.Range("N32").Value = -4675.3335
$ws.Range("H61").Value = 1054.1666
$ws.Range("I61").Value = 716.3158
$ws.Range("J61").Value = 2338
$ws.Range("K61").Value = 716.3158
$ws.Range("L61").Value = 2338
$ws.Range("M61").Value = -504.3158
$ws.Range("N61").Value = -2762
$ws.Range("H74").Value = 1152
$ws.Range("I74").Value = 1004
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 1004
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -130
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 1152
$ws.Range("I77").Value = 1004
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 5020
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -652
$ws.Range("N77").Value = -15236
$ws.Range("H110").Value = 2458.923
$ws.Range("I110").Value = 2701.4285
$ws.Range("J110").Value = 2176
$ws.Range("K110").Value = 2701.4285
$ws.Range("L110").Value = 2176
$ws.Range("M110").Value = -656.4285
$ws.Range("N110").Value = -6266
$ws.Range("H132").Value = 3813.9092
$ws.Range("I132").Value = 3485.0977
$ws.Range("J132").Value = 4776.857
$ws.Range("K132").Value = 10455.2931
$ws.Range("L132").Value = 14330.571
$ws.Range("M132").Value = -7925.293099999999
$ws.Range("N132").Value = -19390.571
$ws.Range("H136").Value = 1054.1666
$ws.Range("I136").Value = 716.3158
$ws.Range("J136").Value = 2338
$ws.Range("K136").Value = 2148.9474
$ws.Range("L136").Value = 7014
$ws.Range("M136").Value = 401.0526
$ws.Range("N136").Value = -12114

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19301.107
$ws.Range("I134").Value = 24408.814
$ws.Range("J134").Value = 2406.3845
$ws.Range("K134").Value = 73226.442
$ws.Range("L134").Value = 7219.1535
$ws.Range("M134").Value = -70691.442
$ws.Range("N134").Value = -12289.1535

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1333.25
$ws.Range("I16").Value = 771.1111
$ws.Range("K16").Value = 771.1111
$ws.Range("M16").Value = -484.1111
$ws.Range("H31").Value = 4001929
$ws.Range("I31").Value = 2128.2222
$ws.Range("K31").Value = 2128.2222
$ws.Range("M31").Value = -1833.2222
$ws.Range("H34").Value = 4001929
$ws.Range("I34").Value = 2128.2222
$ws.Range("K34").Value = 2128.2222
$ws.Range("M34").Value = -1926.2222
$ws.Range("H105").Value = 623.43475
$ws.Range("I105").Value = 315.7143
$ws.Range("J105").Value = 1102.1111
$ws.Range("K105").Value = 315.7143
$ws.Range("L105").Value = 1102.1111
$ws.Range("M105").Value = 1431.2857
$ws.Range("N105").Value = -4596.1111
$ws.Range("H107").Value = 476
$ws.Range("I107").Value = 455.94446
$ws.Range("J107").Value = 656.5
$ws.Range("K107").Value = 455.94446
$ws.Range("L107").Value = 656.5
$ws.Range("M107").Value = 1464.05554
$ws.Range("N107").Value = -4496.5
$ws.Range("H113").Value = 1333.25
$ws.Range("I113").Value = 771.1111
$ws.Range("K113").Value = 771.1111
$ws.Range("M113").Value = 1398.8889
$ws.Range("H132").Value = 2773.5667
$ws.Range("I132").Value = 2200.4092
$ws.Range("J132").Value = 4349.75
$ws.Range("K132").Value = 6601.2276
$ws.Range("L132").Value = 13049.25
$ws.Range("M132").Value = -4071.2276
$ws.Range("N132").Value = -18109.25
$ws.Range("H134").Value = 811.8333
$ws.Range("I134").Value = 594.2
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 1782.6
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = 752.3999999999999
$ws.Range("N134").Value = -10770

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 305.95
$ws.Range("I92").Value = 168.33333
$ws.Range("J92").Value = 512.375
$ws.Range("K92").Value = 504.99999
$ws.Range("L92").Value = 1537.125
$ws.Range("M92").Value = 743.00001
$ws.Range("N92").Value = -4033.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1106.561
$ws.Range("I102").Value = 923.2414
$ws.Range("J102").Value = 1549.5834
$ws.Range("K102").Value = 923.2414
$ws.Range("L102").Value = 1549.5834
$ws.Range("M102").Value = 698.7586
$ws.Range("N102").Value = -4793.5834
$ws.Range("H107").Value = 683
$ws.Range("I107").Value = 504.0909
$ws.Range("J107").Value = 1175
$ws.Range("K107").Value = 504.0909
$ws.Range("L107").Value = 1175
$ws.Range("M107").Value = 1415.9091
$ws.Range("N107").Value = -5015
$ws.Range("H113").Value = 22729482
$ws.Range("I113").Value = 83334500
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 83334500
$ws.Range("L113").Value = 2600
$ws.Range("M113").Value = -83332330
$ws.Range("N113").Value = -6940
$ws.Range("H122").Value = 6024.353
$ws.Range("I122").Value = 15950
$ws.Range("J122").Value = 2970.3076
$ws.Range("K122").Value = 47850
$ws.Range("L122").Value = 8910.9228
$ws.Range("M122").Value = -45400
$ws.Range("N122").Value = -13810.9228
$ws.Range("H132").Value = 62973.574
$ws.Range("I132").Value = 92591.59
$ws.Range("J132").Value = 3737.5454
$ws.Range("K132").Value = 277774.77
$ws.Range("L132").Value = 11212.6362
$ws.Range("M132").Value = -275244.77
$ws.Range("N132").Value = -16272.6362

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 997.6799999999999
$ws.Range("I68").Value = 1014.7826
$ws.Range("J68").Value = 801
$ws.Range("K68").Value = 1014.7826
$ws.Range("L68").Value = 801
$ws.Range("M68").Value = -265.7826
$ws.Range("N68").Value = -2299
$ws.Range("H71").Value = 997.6799999999999
$ws.Range("I71").Value = 1014.7826
$ws.Range("J71").Value = 801
$ws.Range("K71").Value = 5073.913
$ws.Range("L71").Value = 4005
$ws.Range("M71").Value = -1329.913
$ws.Range("N71").Value = -11493
$ws.Range("H82").Value = 1538.25
$ws.Range("I82").Value = 4000
$ws.Range("J82").Value = 717.6667
$ws.Range("K82").Value = 4000
$ws.Range("L82").Value = 717.6667
$ws.Range("M82").Value = -3639
$ws.Range("N82").Value = -1439.6667
$ws.Range("H85").Value = 1538.25
$ws.Range("I85").Value = 4000
$ws.Range("J85").Value = 717.6667
$ws.Range("K85").Value = 4000
$ws.Range("L85").Value = 717.6667
$ws.Range("M85").Value = -2752
$ws.Range("N85").Value = -3213.6667
$ws.Range("H93").Value = 1352570.2
$ws.Range("I93").Value = 2080140.1
$ws.Range("J93").Value = 1369
$ws.Range("K93").Value = 2080140.1
$ws.Range("L93").Value = 1369
$ws.Range("M93").Value = -2078892.1
$ws.Range("N93").Value = -3865
$ws.Range("H104").Value = 25925.25
$ws.Range("J104").Value = 25925.25
$ws.Range("L104").Value = 25925.25
$ws.Range("N104").Value = -32913.25
$ws.Range("H122").Value = 2352.3044
$ws.Range("I122").Value = 2540.8
$ws.Range("J122").Value = 2207.3076
$ws.Range("K122").Value = 7622.400000000001
$ws.Range("L122").Value = 6621.9228
$ws.Range("M122").Value = -5172.400000000001
$ws.Range("N122").Value = -11521.9228
$ws.Range("H136").Value = 9011.5
$ws.Range("I136").Value = 13578.223
$ws.Range("J136").Value = 3140
$ws.Range("K136").Value = 40734.669
$ws.Range("L136").Value = 9420
$ws.Range("M136").Value = -38184.669
$ws.Range("N136").Value = -14520

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1372.8
$ws.Range("I107").Value = 1776.8572
$ws.Range("J107").Value = 430
$ws.Range("K107").Value = 5330.571599999999
$ws.Range("L107").Value = 1290
$ws.Range("M107").Value = -3410.571599999999
$ws.Range("N107").Value = -5130
$ws.Range("H122").Value = 2014.2222
$ws.Range("I122").Value = 2014.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6042.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3592.6666
$ws.Range("H126").Value = 1224.75
$ws.Range("I126").Value = 1989.8182
$ws.Range("J126").Value = 729.7059
$ws.Range("K126").Value = 5969.4546
$ws.Range("L126").Value = 2189.1177
$ws.Range("M126").Value = -3499.4546
$ws.Range("N126").Value = -7129.117700000001
$ws.Range("H132").Value = 1242.6097
$ws.Range("I132").Value = 942.64703
$ws.Range("J132").Value = 2699.5715
$ws.Range("K132").Value = 942.64703
$ws.Range("L132").Value = 2699.5715
$ws.Range("M132").Value = -297.9410899999998
$ws.Range("N132").Value = -13158.7145
$ws.Range("N122").ClearContents()
